# CN Format First Commit...
# Insert the new "cn_invoice" and "cn_vehicle" blocks above the existing
# "cn_acknowledgment" block (which shifts from row 49 down to row 63).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push the existing "cn_acknowledgment" block (rows 49-57) down by 14 rows,
# freeing rows 49-62 for the two new blocks.
$ws.Rows.Item(49).Resize(14).Insert()

# ---------------------------------------------------------------------
# Block: cn_invoice (new header at row 49)
# ---------------------------------------------------------------------
$ws.Range("A63:E63").Copy()
$ws.Range("A49:E49").PasteSpecial(-4122)
$ws.Range("A49").Value() = "cn_invoice"

$ws.Range("A64").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("C64").Copy()
$ws.Range("C50").PasteSpecial(-4122)
$ws.Range("E64").Copy()
$ws.Range("E50").PasteSpecial(-4122)
$ws.Range("A50").Value() = "Field Name"
$ws.Range("C50").Value() = "Field Name"
$ws.Range("E50").Value() = "Field Name"

$ws.Range("A51").Value() = "cn_id"

$ws.Range("A52").Value() = "cn_invoice_no"
$ws.Range("C52").Value() = "cn_invoice_total"
$ws.Range("E52").Value() = "cn_invoice_tax"

$ws.Range("A53").Value() = "cn_invoice_date"
$ws.Range("C53").Value() = "cn_invoice_gstin"
$ws.Range("E53").Value() = "cn_invoice_quantity"

$ws.Range("A54").Value() = "cn_invoice_taxable"
$ws.Range("C54").Value() = "cn_invoice_weight"
$ws.Range("E54").Value() = "cn_invoice_mop"

# row 55 intentionally left blank (spacer row, like the other blocks)

# ---------------------------------------------------------------------
# Block: cn_vehicle (new header at row 56)
# ---------------------------------------------------------------------
$ws.Range("A63:E63").Copy()
$ws.Range("A56:E56").PasteSpecial(-4122)
$ws.Range("A56").Value() = "cn_vehicle"
$ws.Range("F56").ClearContents()

$ws.Range("A64").Copy()
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("C64").Copy()
$ws.Range("C57").PasteSpecial(-4122)
$ws.Range("E64").Copy()
$ws.Range("E57").PasteSpecial(-4122)
$ws.Range("A57").Value() = "Field Name"
$ws.Range("C57").Value() = "Field Name"
$ws.Range("E57").Value() = "Field Name"

$ws.Range("A58").Value() = "cn_id"

$ws.Range("A59").Value() = "cn_vehicle_no."
$ws.Range("C59").Value() = "cn_vehicle"
$ws.Range("E59").Value() = "cn_vehicle"

$ws.Range("A60").Value() = "cn_vehicle_type"
$ws.Range("C60").Value() = "cn_vehicle"
$ws.Range("E60").Value() = "cn_vehicle"

$ws.Range("A61").Value() = "cn_vehicle"
$ws.Range("C61").Value() = "cn_vehicle"
$ws.Range("E61").Value() = "cn_vehicle"

# row 62 intentionally left blank (spacer row, like the other blocks)

# ---------------------------------------------------------------------
# Register the merged header ranges for the two new blocks.
# ---------------------------------------------------------------------
$ws.Range("A49:E49").Merge()
$ws.Range("A56:E56").Merge()

# ---------------------------------------------------------------------
# Restore view state (best effort): scroll + selection.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow() = 44
$excel.ActiveWindow.ScrollColumn() = 1
$ws.Range("M61").Select()
